# Natmi following Dr Hou advice
# Re-run of the LR-pair analysis for Pdgfc-Pdgfra adds an "ECs" cluster
# alongside the existing "FAPs"/"sCs" clusters, expanding the 2-row result
# table (2 sending x 1 target cluster) into the full 3x3 sending/target grid.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Build the 9-row x 20-column (A:T) result block as a single 2D array so it
# can be written in one Range.Value assignment (rows 2-10).
$arr = New-Object 'object[,]' 9,20

# Row 2: ECs -> ECs (via Pdgfc/Pdgfra)
$arr[0,0] = "ECs"
$arr[0,1] = "Pdgfc"
$arr[0,2] = "Pdgfra"
$arr[0,3] = "ECs"
$arr[0,4] = 2
$arr[0,5] = 0.6666666666666666
$arr[0,6] = 0.2015403333333333
$arr[0,7] = 0.604621
$arr[0,8] = 0.05371480119195454
$arr[0,9] = 0.05371480119195454
$arr[0,10] = 2
$arr[0,11] = 0.6666666666666666
$arr[0,12] = 1.155747666666667
$arr[0,13] = 3.467243
$arr[0,14] = 0.004246591903937912
$arr[0,15] = 0.004246591903937912
$arr[0,16] = 0.2329297699892222
$arr[0,17] = 2.096367929903
$arr[0,18] = 0.0002281048398633887
$arr[0,19] = 0.0002281048398633887

# Row 3: ECs -> FAPs (via Pdgfc/Pdgfra)
$arr[1,0] = "ECs"
$arr[1,1] = "Pdgfc"
$arr[1,2] = "Pdgfra"
$arr[1,3] = "FAPs"
$arr[1,4] = 2
$arr[1,5] = 0.6666666666666666
$arr[1,6] = 0.2015403333333333
$arr[1,7] = 0.604621
$arr[1,8] = 0.05371480119195454
$arr[1,9] = 0.05371480119195454
$arr[1,10] = 3
$arr[1,11] = 1
$arr[1,12] = 270.7963256666667
$arr[1,13] = 812.3889770000001
$arr[1,14] = 0.9949935590256014
$arr[1,15] = 0.9949935590256014
$arr[1,16] = 54.57638174030189
$arr[1,17] = 491.187435662717
$arr[1,18] = 0.05344588121033546
$arr[1,19] = 0.05344588121033546

# Row 4: ECs -> sCs (via Pdgfc/Pdgfra)
$arr[2,0] = "ECs"
$arr[2,1] = "Pdgfc"
$arr[2,2] = "Pdgfra"
$arr[2,3] = "sCs"
$arr[2,4] = 2
$arr[2,5] = 0.6666666666666666
$arr[2,6] = 0.2015403333333333
$arr[2,7] = 0.604621
$arr[2,8] = 0.05371480119195454
$arr[2,9] = 0.05371480119195454
$arr[2,10] = 3
$arr[2,11] = 1
$arr[2,12] = 0.2067996666666667
$arr[2,13] = 0.620399
$arr[2,14] = 0.0007598490704606447
$arr[2,15] = 0.0007598490704606446
$arr[2,16] = 0.04167847375322222
$arr[2,17] = 0.375106263779
$arr[2,18] = 0.00004081514175568498
$arr[2,19] = 0.00004081514175568498

# Row 5: FAPs -> ECs (via Pdgfc/Pdgfra)
$arr[3,0] = "FAPs"
$arr[3,1] = "Pdgfc"
$arr[3,2] = "Pdgfra"
$arr[3,3] = "ECs"
$arr[3,4] = 2
$arr[3,5] = 0.6666666666666666
$arr[3,6] = 1.961473666666667
$arr[3,7] = 5.884421
$arr[3,8] = 0.5227746044956465
$arr[3,9] = 0.5227746044956465
$arr[3,10] = 2
$arr[3,11] = 0.6666666666666666
$arr[3,12] = 1.155747666666667
$arr[3,13] = 3.467243
$arr[3,14] = 0.004246591903937912
$arr[3,15] = 0.004246591903937912
$arr[3,16] = 2.266968613478111
$arr[3,17] = 20.402717521303
$arr[3,18] = 0.002220010403035556
$arr[3,19] = 0.002220010403035556

# Row 6: FAPs -> FAPs (via Pdgfc/Pdgfra)
$arr[4,0] = "FAPs"
$arr[4,1] = "Pdgfc"
$arr[4,2] = "Pdgfra"
$arr[4,3] = "FAPs"
$arr[4,4] = 2
$arr[4,5] = 0.6666666666666666
$arr[4,6] = 1.961473666666667
$arr[4,7] = 5.884421
$arr[4,8] = 0.5227746044956465
$arr[4,9] = 0.5227746044956465
$arr[4,10] = 3
$arr[4,11] = 1
$arr[4,12] = 270.7963256666667
$arr[4,13] = 812.3889770000001
$arr[4,14] = 0.9949935590256014
$arr[4,15] = 0.9949935590256014
$arr[4,16] = 531.1598618252575
$arr[4,17] = 4780.438756427317
$arr[4,18] = 0.5201573642953244
$arr[4,19] = 0.5201573642953244

# Row 7: FAPs -> sCs (via Pdgfc/Pdgfra)
$arr[5,0] = "FAPs"
$arr[5,1] = "Pdgfc"
$arr[5,2] = "Pdgfra"
$arr[5,3] = "sCs"
$arr[5,4] = 2
$arr[5,5] = 0.6666666666666666
$arr[5,6] = 1.961473666666667
$arr[5,7] = 5.884421
$arr[5,8] = 0.5227746044956465
$arr[5,9] = 0.5227746044956465
$arr[5,10] = 3
$arr[5,11] = 1
$arr[5,12] = 0.2067996666666667
$arr[5,13] = 0.620399
$arr[5,14] = 0.0007598490704606447
$arr[5,15] = 0.0007598490704606446
$arr[5,16] = 0.4056321004421111
$arr[5,17] = 3.650688903979
$arr[5,18] = 0.0003972297972864482
$arr[5,19] = 0.0003972297972864482

# Row 8: sCs -> ECs (via Pdgfc/Pdgfra)
$arr[6,0] = "sCs"
$arr[6,1] = "Pdgfc"
$arr[6,2] = "Pdgfra"
$arr[6,3] = "ECs"
$arr[6,4] = 3
$arr[6,5] = 1
$arr[6,6] = 1.589030666666667
$arr[6,7] = 4.767092
$arr[6,8] = 0.423510594312399
$arr[6,9] = 0.4235105943123989
$arr[6,10] = 2
$arr[6,11] = 0.6666666666666666
$arr[6,12] = 1.155747666666667
$arr[6,13] = 3.467243
$arr[6,14] = 0.004246591903937912
$arr[6,15] = 0.004246591903937912
$arr[6,16] = 1.836518485261778
$arr[6,17] = 16.528666367356
$arr[6,18] = 0.001798476661038967
$arr[6,19] = 0.001798476661038967

# Row 9: sCs -> FAPs (via Pdgfc/Pdgfra)
$arr[7,0] = "sCs"
$arr[7,1] = "Pdgfc"
$arr[7,2] = "Pdgfra"
$arr[7,3] = "FAPs"
$arr[7,4] = 3
$arr[7,5] = 1
$arr[7,6] = 1.589030666666667
$arr[7,7] = 4.767092
$arr[7,8] = 0.423510594312399
$arr[7,9] = 0.4235105943123989
$arr[7,10] = 3
$arr[7,11] = 1
$arr[7,12] = 270.7963256666667
$arr[7,13] = 812.3889770000001
$arr[7,14] = 0.9949935590256014
$arr[7,15] = 0.9949935590256014
$arr[7,16] = 430.3036659049872
$arr[7,17] = 3872.732993144884
$arr[7,18] = 0.4213903135199414
$arr[7,19] = 0.4213903135199414

# Row 10: sCs -> sCs (via Pdgfc/Pdgfra)
$arr[8,0] = "sCs"
$arr[8,1] = "Pdgfc"
$arr[8,2] = "Pdgfra"
$arr[8,3] = "sCs"
$arr[8,4] = 3
$arr[8,5] = 1
$arr[8,6] = 1.589030666666667
$arr[8,7] = 4.767092
$arr[8,8] = 0.423510594312399
$arr[8,9] = 0.4235105943123989
$arr[8,10] = 3
$arr[8,11] = 1
$arr[8,12] = 0.2067996666666667
$arr[8,13] = 0.620399
$arr[8,14] = 0.0007598490704606447
$arr[8,15] = 0.0007598490704606446
$arr[8,16] = 0.3286110121897778
$arr[8,17] = 2.957499109708
$arr[8,18] = 0.0003218041314185116
$arr[8,19] = 0.0003218041314185115

$ws.Range("A2:T10").Value = $arr

Write-Output "Wrote $([int]1+9 -1) data rows (A2:T10)"
